$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# B11 held the shared string "R40" (row 11 of the "Rules" table). The rule's
# upper bound changed, so the label becomes the literal text "1". A leading
# apostrophe forces Excel to store it as text (shared string) rather than
# coercing the numeric-looking literal to a number.
$ws.Range("B11").Value = "'1"
